$d = $word.ActiveDocument

# --- Add the three new character styles (matching the diff added to styles.xml) ---

$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1

# --- Apply GaNStyle to every "2022 Campaign Dates ..." run (4 occurrences) ---

$range = $d.Content
while ($range.Find.Execute(" 2022 Campaign Dates that use Bootes constellation: May 14-23, June 13-22, July 12-21", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $range.Style = "GaNStyle"
    $range.Collapse(0)
}

# --- Apply GaNParagraph to the "You are participating..." run ---

$range = $d.Content
if ($range.Find.Execute("You are participating in a global campaign to observe and record the faintest stars visible as a means of measuring light pollution in a given location. By locating and observing the constellation Bootes constellation in the night sky and comparing it to stellar charts, people from around the world will learn how the lights in their community contribute to light pollution. Your contributions to the online database will document the visible nighttime sky.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $range.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)." run ---

$range = $d.Content
if ($range.Find.Execute("(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $range.Style = "GaNLinks"
}
